$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.04352703914852966
$ws.Range("C2").Value = -0.1859241290933605
$ws.Range("B3").Value = -0.02565587424673664
$ws.Range("C3").Value = -0.347112309470367
$ws.Range("B4").Value = 0.09130917296212129
$ws.Range("C4").Value = 0.5205335838175501
$ws.Range("B5").Value = 0.1407264577336982
$ws.Range("C5").Value = -0.3562054692811555
$ws.Range("B6").Value = 0.1365801548787149
$ws.Range("C6").Value = 0.03191918419370631
$ws.Range("B7").Value = 0.4256091460202764
$ws.Range("C7").Value = 0.1326234609464499
$ws.Range("B8").Value = 0.387497880336328
$ws.Range("C8").Value = -0.3793809547011148
$ws.Range("B9").Value = 0.07970100427750651
$ws.Range("C9").Value = 0.2703682930351588
$ws.Range("B10").Value = 0.5529856218129925
$ws.Range("C10").Value = -0.0180458906639075
$ws.Range("B11").Value = 0.1573638192463827
$ws.Range("C11").Value = 0.06349312328929377
$ws.Range("B12").Value = -0.3638743103995788
$ws.Range("C12").Value = 0.007665347098700104
$ws.Range("B13").Value = -0.3532615717166662
$ws.Range("C13").Value = -0.2622057855656936
$ws.Range("B14").Value = -0.1406593506252099
$ws.Range("C14").Value = 0.2785728733302939
$ws.Range("B15").Value = -0.07374420783335736
$ws.Range("C15").Value = -0.2200253420549833
$ws.Range("B16").Value = 0.0029183561513963
$ws.Range("C16").Value = -0.1042442078322221
$ws.Range("B17").Value = -0.00390275612779649
$ws.Range("C17").Value = 0.03771552799532983
